$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage so that
# numeric-looking strings (e.g. "1.009", "27.013.04") are not silently
# reinterpreted by Excel as numbers. Restoring the cell Style afterwards
# keeps the original (unstyled) cell formatting intact.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '27.013.04'
Set-TextValue 'E2' '  -1.25%  '

# Row 3
Set-TextValue 'D3' '1.822.36'
Set-TextValue 'E3' '  -0.62%  '

# Row 4
Set-TextValue 'D4' '1.009'
Set-TextValue 'E4' '  -0.29%  '

# Row 5
Set-TextValue 'D5' '309.57'
Set-TextValue 'E5' '  -1.64%  '

# Row 6
Set-TextValue 'D6' '1.008'
Set-TextValue 'E6' '  -0.28%  '

# Row 7
Set-TextValue 'D7' '0.4656'
Set-TextValue 'E7' '  -2.02%  '

# Row 8
Set-TextValue 'D8' '0.3659'
Set-TextValue 'E8' '  -0.80%  '

# Row 9
Set-TextValue 'D9' '0.07234'
Set-TextValue 'E9' '  -3.00%  '

# Row 10
Set-TextValue 'D10' '0.8599'
Set-TextValue 'E10' '  -2.93%  '

# Row 11
Set-TextValue 'D11' '19.84'
Set-TextValue 'E11' '  -2.92%  '

# Row 12
Set-TextValue 'D12' '0.07546'
Set-TextValue 'E12' '  +3.06%  '

# Row 13
Set-TextValue 'D13' '1.810.35'
Set-TextValue 'E13' '  -4.00%  '

# Row 14
Set-TextValue 'D14' '5.332'
Set-TextValue 'E14' '  -2.04%  '

# Row 15
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '6.496'
Set-TextValue 'E15' '  -1.25%  '

# Row 16
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D16' '91.59'
Set-TextValue 'E16' '  -1.70%  '

# Row 17
Set-TextValue 'D17' '1.009'
Set-TextValue 'E17' '  -0.12%  '

# Row 18
Set-TextValue 'E18' '  -1.91%  '

# Row 19
Set-TextValue 'D19' '1.007'
Set-TextValue 'E19' '  -0.39%  '

# Row 20
Set-TextValue 'D20' '14.48'
Set-TextValue 'E20' '  -2.19%  '

# Row 21
Set-TextValue 'D21' '26.883.96'
Set-TextValue 'E21' '  -2.61%  '

# Row 22
Set-TextValue 'D22' '5.146'
Set-TextValue 'E22' '  -2.73%  '

# Row 23
Set-TextValue 'D23' '10.52'
Set-TextValue 'E23' '  -1.30%  '

# Row 24
Set-TextValue 'D24' '1.984.93'
Set-TextValue 'E24' '  -5.37%  '

# Row 25
Set-TextValue 'D25' '151.44'
Set-TextValue 'E25' '  -0.25%  '

# Row 26
Set-TextValue 'D26' '1.841'
Set-TextValue 'E26' '  -2.68%  '

# Row 27
Set-TextValue 'D27' '18.12'
Set-TextValue 'E27' '  -2.89%  '

# Row 28
Set-TextValue 'D28' '2.052'
Set-TextValue 'E28' '  -4.04%  '

# Row 29
Set-TextValue 'D29' '5.109'
Set-TextValue 'E29' '  -2.50%  '

# Row 30
Set-TextValue 'D30' '115.26'
Set-TextValue 'E30' '  -1.82%  '

# Row 31
Set-TextValue 'D31' '0.08859'
Set-TextValue 'E31' '  -1.55%  '

# Row 32
Set-TextValue 'D32' '2.957'
Set-TextValue 'E32' '  +0.51%  '

# Row 33
Set-TextValue 'D33' '4.426'
Set-TextValue 'E33' '  -2.73%  '

# Row 34
Set-TextValue 'D34' '1.131'
Set-TextValue 'E34' '  -3.96%  '

# Row 35
Set-TextValue 'D35' '0.7177'
Set-TextValue 'E35' '  -4.69%  '

# Row 36
Set-TextValue 'D36' '1.078'
Set-TextValue 'E36' '  -2.17%  '

# Row 37
Set-TextValue 'D37' '0.05261'
Set-TextValue 'E37' '  -1.62%  '

# Row 38
Set-TextValue 'D38' '0.01920'
Set-TextValue 'E38' '  -1.89%  '

# Row 39
Set-TextValue 'D39' '2.399'
Set-TextValue 'E39' '  +0.23%  '

# Row 40
Set-TextValue 'D40' '2.927'
Set-TextValue 'E40' '  -1.78%  '

# Row 41
Set-TextValue 'D41' '7.155'
Set-TextValue 'E41' '  -1.76%  '

# Row 42
Set-TextValue 'D42' '0.5156'
Set-TextValue 'E42' '  -3.02%  '

# Row 43
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D43' '0.8596'
Set-TextValue 'E43' '  -14.97%  '

# Row 44
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D44' '0.1625'
Set-TextValue 'E44' '  -2.12%  '

# Row 45
Set-TextValue 'D45' '8.164'
Set-TextValue 'E45' '  -3.69%  '

# Row 46
Set-TextValue 'D46' '0.4803'
Set-TextValue 'E46' '  -2.26%  '

# Row 47
Set-TextValue 'D47' '1.008'
Set-TextValue 'E47' '  -0.33%  '

# Row 48
Set-TextValue 'E48' '  -3.79%  '

# Row 49
Set-TextValue 'D49' '102.62'
Set-TextValue 'E49' '  -2.36%  '

# Row 51
Set-TextValue 'D51' '0.06245'
Set-TextValue 'E51' '  -0.87%  '
